$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Add the new "2023" sheet right after "2022", mirroring the "Laporan Obat
# UGD" report for the new year.
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "2023"

# Bring over the header row (Bulan / BPJS / Umum / Total) and the month
# labels in column A so the new sheet reuses the same shared strings and
# header styling (bottom border) as the 2022 sheet.
$ws1.Range("A1:A13").Copy($ws2.Range("A1:A13"))
$ws1.Range("B1:D1").Copy($ws2.Range("B1:D1"))

# Totals column - same SUM formula as 2022, written as one range so Excel
# stores it as a shared formula like the original sheet.
$ws2.Range("D2:D13").Formula = "=SUM(B2:C2)"

# Update the 2022 sheet's remembered selection first...
$ws1.Range("A1:E15").Select() | Out-Null

# ...then make the new sheet the active tab, with D27 selected, so it is
# the one left on-screen/selected when the workbook is saved.
$ws2.Activate()
$ws2.Range("D27").Select() | Out-Null
